$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new year column K with header 2022 and data values, matching style of column J
$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 26.495524312074597
$ws.Range("K6").Value = 59.383769502755833
$ws.Range("K7").Value = 38.32334404557426
$ws.Range("K8").Value = 48.136790950525594
$ws.Range("K9").Value = 46.63213064070051
$ws.Range("K10").Value = 32.657429481680126
$ws.Range("K11").Value = 31.457245964894081
$ws.Range("K12").Value = 22.734405597714229
$ws.Range("K13").Value = -0.19691879995369213
$ws.Range("K14").Value = 33.158040409631916

# Copy style/formatting from column J to column K for rows 4-14
$ws.Range("J4:J14").Copy()
$ws.Range("K4:K14").PasteSpecial(-4122)  # xlPasteFormats

# Update the selection to match the target state
$ws.Range("M7").Select()
